$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C2:C9) from 45175 to 45183 for all data rows
$ws.Range("C2:C9").Value = 45183
